# Auto-generated edit script applying the "Updated symbol list" commit.
# All numeric-looking text cells (columns D and G) are written through
# a text NumberFormat so Excel keeps them as strings (matching the
# original inlineStr cells) instead of silently coercing to numbers;
# the style is then reset to "Normal" so no stray formatting remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("Hora") : every data row goes from 15 -> 16 ---
$gRange = $ws.Range("G2:G51")
$gRange.NumberFormat = "@"
$gRange.Value = "16"
$gRange.Style = "Normal"

# --- Column D ("Price") : updated quote values ---
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"
$ws.Range("D2").Value = "250.53"
$ws.Range("D3").Value = "22.87"
$ws.Range("D4").Value = "5.452"
$ws.Range("D5").Value = "0.05627"
$ws.Range("D6").Value = "3.431"
$ws.Range("D7").Value = "6.403"
$ws.Range("D8").Value = "0.8184"
$ws.Range("D9").Value = "0.9244"
$ws.Range("D10").Value = "0.1434"
$ws.Range("D11").Value = "0.07514"
$ws.Range("D12").Value = "0.03156"
$ws.Range("D13").Value = "0.03088"
$ws.Range("D14").Value = "0.09335"
$ws.Range("D15").Value = "3.562"
$ws.Range("D16").Value = "0.001615"
$ws.Range("D17").Value = "0.04742"
$ws.Range("D18").Value = "0.006435"
$ws.Range("D19").Value = "0.004991"
$ws.Range("D20").Value = "0.001034"
$ws.Range("D21").Value = "0.0001502"
$ws.Range("D22").Value = "3.720"
$ws.Range("D23").Value = "2.176"
$ws.Range("D24").Value = "0.01151"
$ws.Range("D25").Value = "0.3308"
$ws.Range("D26").Value = "0.1275"
$ws.Range("D28").Value = "0.0003004"
$ws.Range("D40").Value = "0.04013"
$ws.Range("D41").Value = "0.006838"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").Value = "0.002766"
$ws.Range("D44").Value = "0.007878"
$ws.Range("D45").Value = "0.00005579"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D48").Value = "0.5006"
$ws.Range("D49").Value = "0.2256"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D51").Value = "0.01011"
$dRange.Style = "Normal"

# --- Columns B/C/E : coin re-rank / rename / link updates ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"
$ws.Range("E49").Value = "48BOLOBOLO"
